# Swap the order of recorders in the "Recorded By" column (column G) so that
# any cell listing both "dnasr281@gmail.com" and "System" shows "System" first,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Cells that only contain a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
